$wb = $excel.ActiveWorkbook
$phases = $wb.Worksheets.Item("Phases")
$events = $wb.Worksheets.Item("Events")

# ============================================================
# Events sheet: new "SortOrder" column (J) + new event row (3)
# ============================================================
$events.Range("C3").Value = "Update/Tweet: PCP Coming Soon"
$events.Range("D2").Copy()
$events.Range("D3").PasteSpecial(-4122)
$events.Range("D3").Value = "Notification Review"
$events.Range("E3").Value = "Communications"
$events.Range("F3").Value = "Calendar"

$events.Range("A3").Value = 2
$events.Range("B3").Value = 1
$events.Range("G3").Value = 7
$events.Range("H3").Value = 1
$events.Range("I3").Value = $true

$events.Range("J1").Value = "SortOrder"
$events.Range("J2").Value = 1
$events.Range("J3").Value = 2

$events.Range("C3").Interior.Color = 16764108
$events.Range("C3").VerticalAlignment = -4108

# ============================================================
# Phases sheet: new phase row (3), then a new "Legislated"
# column inserted before "Color"
# ============================================================
$phases.Range("B3").Value = "Early Engagement"
$phases.Range("C3").Value = "Assessment"
$phases.Range("D3").Value = "EA Act (2018)"
$phases.Range("A3").Value = 2
$phases.Range("E3").Value = 60
$phases.Range("F3").Value = "#54858d"
$phases.Range("G3").Value = 2

$phases.Columns.Item(6).Insert()
$phases.Range("F1").Value = "Legislated"
$phases.Range("F2").Value = $true
$phases.Range("F3").Value = $true
$phases.Columns.Item(6).ColumnWidth = 12.666666666666666

# ============================================================
# View / selection state
# ============================================================
$phases.Range("F4").Select()
$events.Activate()
$events.Range("C9").Select()
